$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1: "My 10 years ... Hempel's GrowHub fulfill its mission ..."
#            -> "My Ph.D. in ... Nilfisk innovate in a sustainable way ..."
# ---------------------------------------------------------------------------

$r1a = $d.Content
$r1a.Find.Execute(
    "My 10 years of experience in developing and managing international research projects both as a leader and collaborator, along with my multidisciplinary background, provide me with the skills to help Hempel's ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "My Ph.D. in Materials Science and Engineering and more than 2 years as a postdoctoral fellow in computational solid and fluid mechanics and advanced fracture mechanics, along with my multidisciplinary engineering background, provide me with the skills to help ",
    2) | Out-Null

$r1b = $d.Content
$r1b.Find.Execute("GrowHub", $false, $false, $false, $false, $false, $true, 1, $false, "Nilfisk", 2) | Out-Null

$r1c = $d.Content
$r1c.Find.Execute(
    " fulfill its mission in an Innovation Manager capacity.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " innovate in a sustainable way in a R&D Mechanical Engineer capacity.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 2: "I offer proficiency in managing international R&D projects..."
#            -> "During my professional journey I contributed to projects..."
# ---------------------------------------------------------------------------

$r2 = $d.Content
$r2.Find.Execute(
    "I offer proficiency in managing international R&D projects with uncertain boundaries, evolving requirements, and multiple stakeholders from ideation to exploitation, as well as experience in research funding and grant writing. I am expert in 3D modeling and research software development, and several digital tools and programming languages as outlined in my CV.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "During my professional journey I contributed to projects in several fields, from ankle biomechanics to the design of nano-sized electromagnetic sensors, from modeling and prediction of damage in fiber-reinforced composites to multi-scale modeling of wood, from coupled fluid-structure interaction to large displacement analysis of cracking in hydrogels. I offer proficiency in several CAD and CAE tools, as well as computational methods of solid mechanics (FEM, BEM), fluid mechanics (LBM, FVM), fracture and damage mechanics (CZM, VCCT, J-integral, interaction integrals), mesh generation and computational geometry (Delaunay triangulation, transfinite interpolation, elliptic/parabolic/hyperbolic mesh smoothing). I have multiple years of experience in mechanical testing and I am expert in several programming languages, as outlined in my CV.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 3: "... in service of GrowHub's mission ..."
#            -> "... in service of Nilfisk's mission ..."
# ---------------------------------------------------------------------------

$r3 = $d.Content
$r3.Find.Execute("GrowHub's", $false, $false, $false, $false, $false, $true, 1, $false, "Nilfisk's", 2) | Out-Null

Write-Host "Done"
